$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9524916
$ws.Range("I40").Value = 13334330
$ws.Range("J40").Value = 1380
$ws.Range("K40").Value = 13334330
$ws.Range("L40").Value = 1380
$ws.Range("M40").Value = -13334155
$ws.Range("N40").Value = -1730
$ws.Range("H74").Value = 4000
$ws.Range("I74").Value = 3866.6667
$ws.Range("J74").Value = 4266.6665
$ws.Range("K74").Value = 3866.6667
$ws.Range("L74").Value = 4266.6665
$ws.Range("M74").Value = -2930.6667
$ws.Range("N74").Value = -6138.6665
$ws.Range("H76").Value = 3075
$ws.Range("I76").Value = 3060
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 3060
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -2745
$ws.Range("N76").Value = -3830
$ws.Range("H77").Value = 4000
$ws.Range("I77").Value = 3866.6667
$ws.Range("J77").Value = 4266.6665
$ws.Range("K77").Value = 19333.3335
$ws.Range("L77").Value = 21333.3325
$ws.Range("M77").Value = -14653.3335
$ws.Range("N77").Value = -30693.3325
$ws.Range("H79").Value = 3075
$ws.Range("I79").Value = 3060
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 3060
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -1968
$ws.Range("N79").Value = -5384
$ws.Range("H118").Value = 951.11536
$ws.Range("I118").Value = 645
$ws.Range("J118").Value = 1440.9
$ws.Range("K118").Value = 1935
$ws.Range("L118").Value = 4322.700000000001
$ws.Range("M118").Value = -278
$ws.Range("N118").Value = -7636.700000000001
$ws.Range("H128").Value = 28515.385
$ws.Range("J128").Value = 28515.385
$ws.Range("L128").Value = 28515.385
$ws.Range("N128").Value = -38475.38499999999
$ws.Range("H129").Value = 2157
$ws.Range("I129").Value = 3100
$ws.Range("J129").Value = 2071.2727
$ws.Range("K129").Value = 9300
$ws.Range("L129").Value = 6213.8181
$ws.Range("M129").Value = -4300
$ws.Range("N129").Value = -16213.8181
$ws.Range("H132").Value = 2516.8223
$ws.Range("I132").Value = 2482.721
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 7448.163
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -4918.163
$ws.Range("N132").Value = -14810
$ws.Range("H138").Value = 2349.9023
$ws.Range("I138").Value = 2257.6538
$ws.Range("J138").Value = 2509.8
$ws.Range("K138").Value = 6772.9614
$ws.Range("L138").Value = 7529.400000000001
$ws.Range("M138").Value = -1632.9614
$ws.Range("N138").Value = -17809.4
$ws.Range("H141").Value = 1364.1177
$ws.Range("I141").Value = 1008.1818
$ws.Range("J141").Value = 2016.6666
$ws.Range("K141").Value = 3024.5454
$ws.Range("L141").Value = 6049.9998
$ws.Range("M141").Value = 2155.4546
$ws.Range("N141").Value = -16409.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2263058
$ws.Range("I2").Value = 586.375
$ws.Range("J2").Value = 5883013
$ws.Range("K2").Value = 586.375
$ws.Range("L2").Value = 5883013
$ws.Range("M2").Value = -473.375
$ws.Range("N2").Value = -5883239
$ws.Range("H32").Value = 7927.8027
$ws.Range("I32").Value = 6129.271
$ws.Range("J32").Value = 14169.765
$ws.Range("K32").Value = 6129.271
$ws.Range("L32").Value = 14169.765
$ws.Range("M32").Value = -5842.271
$ws.Range("N32").Value = -14743.765
$ws.Range("H116").Value = 2263058
$ws.Range("I116").Value = 586.375
$ws.Range("J116").Value = 5883013
$ws.Range("K116").Value = 586.375
$ws.Range("L116").Value = 5883013
$ws.Range("M116").Value = 1707.625
$ws.Range("N116").Value = -5887601
$ws.Range("H122").Value = 1211.4584
$ws.Range("I122").Value = 1235.7142
$ws.Range("K122").Value = 3707.1426
$ws.Range("M122").Value = -1257.1426

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2263058
$ws.Range("I3").Value = 586.375
$ws.Range("J3").Value = 5883013
$ws.Range("K3").Value = 586.375
$ws.Range("L3").Value = 5883013
$ws.Range("M3").Value = -472.375
$ws.Range("N3").Value = -5883241

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1307
$ws.Range("I99").Value = 1295.1111
$ws.Range("J99").Value = 1414
$ws.Range("K99").Value = 1295.1111
$ws.Range("L99").Value = 1414
$ws.Range("M99").Value = 202.8888999999999
$ws.Range("N99").Value = -4410
$ws.Range("H107").Value = 1001.8929
$ws.Range("I107").Value = 1269.8572
$ws.Range("J107").Value = 198
$ws.Range("K107").Value = 1269.8572
$ws.Range("L107").Value = 198
$ws.Range("M107").Value = 650.1428000000001
$ws.Range("N107").Value = -4038
$ws.Range("H126").Value = 1307
$ws.Range("I126").Value = 1295.1111
$ws.Range("J126").Value = 1414
$ws.Range("K126").Value = 3885.3333
$ws.Range("L126").Value = 4242
$ws.Range("M126").Value = -1415.3333
$ws.Range("N126").Value = -9182
$ws.Range("H132").Value = 3156.2144
$ws.Range("I132").Value = 2845.7144
$ws.Range("K132").Value = 8537.143199999999
$ws.Range("M132").Value = -6007.143199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 468.0909
$ws.Range("I18").Value = 372.22223
$ws.Range("J18").Value = 899.5
$ws.Range("K18").Value = 1116.66669
$ws.Range("L18").Value = 2698.5
$ws.Range("M18").Value = -947.66669
$ws.Range("N18").Value = -3036.5
$ws.Range("H113").Value = 705.5714
$ws.Range("I113").Value = 681
$ws.Range("J113").Value = 724
$ws.Range("K113").Value = 2043
$ws.Range("L113").Value = 2172
$ws.Range("M113").Value = 127
$ws.Range("N113").Value = -6512
$ws.Range("H114").Value = 1664.2
$ws.Range("I114").Value = 478.4
$ws.Range("J114").Value = 2850
$ws.Range("K114").Value = 1435.2
$ws.Range("L114").Value = 8550
$ws.Range("M114").Value = 1818.8
$ws.Range("N114").Value = -15058
$ws.Range("H126").Value = 6943.3335
$ws.Range("I126").Value = 830
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 2490
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = 2450
$ws.Range("N126").Value = -39880
$ws.Range("H130").Value = 2225.5557
$ws.Range("I130").Value = 1876.6666
$ws.Range("J130").Value = 2400
$ws.Range("K130").Value = 5629.9998
$ws.Range("L130").Value = 7200
$ws.Range("M130").Value = -609.9997999999996
$ws.Range("N130").Value = -17240
$ws.Range("H137").Value = 55174.31
$ws.Range("I137").Value = 2476.923
$ws.Range("J137").Value = 81523
$ws.Range("K137").Value = 7430.768999999999
$ws.Range("L137").Value = 244569
$ws.Range("M137").Value = -2330.768999999999
$ws.Range("N137").Value = -254769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4528.0625
$ws.Range("I70").Value = 4650
$ws.Range("J70").Value = 3999.6667
$ws.Range("K70").Value = 4650
$ws.Range("L70").Value = 3999.6667
$ws.Range("M70").Value = -4380
$ws.Range("N70").Value = -4539.6667
$ws.Range("H73").Value = 4528.0625
$ws.Range("I73").Value = 4650
$ws.Range("J73").Value = 3999.6667
$ws.Range("K73").Value = 4650
$ws.Range("L73").Value = 3999.6667
$ws.Range("M73").Value = -3714
$ws.Range("N73").Value = -5871.6667
$ws.Range("H107").Value = 554.95654
$ws.Range("I107").Value = 351.91666
$ws.Range("J107").Value = 776.4545000000001
$ws.Range("K107").Value = 351.91666
$ws.Range("L107").Value = 776.4545000000001
$ws.Range("M107").Value = 1568.08334
$ws.Range("N107").Value = -4616.4545
$ws.Range("H132").Value = 1607.4
$ws.Range("I132").Value = 1139.3889
$ws.Range("K132").Value = 3418.1667
$ws.Range("M132").Value = -888.1666999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1578.8
$ws.Range("I7").Value = 1473.5
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 1473.5
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -1361.5
$ws.Range("N7").Value = -2224
$ws.Range("H61").Value = 1998.2858
$ws.Range("I61").Value = 1397.6
$ws.Range("K61").Value = 1397.6
$ws.Range("M61").Value = -1195.6
$ws.Range("H113").Value = 1998.2858
$ws.Range("I113").Value = 1397.6
$ws.Range("K113").Value = 1397.6
$ws.Range("M113").Value = 772.4000000000001
$ws.Range("H126").Value = 1578.8
$ws.Range("I126").Value = 1473.5
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 4420.5
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1950.5
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 15600.223
$ws.Range("I132").Value = 100004
$ws.Range("J132").Value = 5049.75
$ws.Range("K132").Value = 300012
$ws.Range("L132").Value = 15149.25
$ws.Range("M132").Value = -297482
$ws.Range("N132").Value = -20209.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 454.58066
$ws.Range("I113").Value = 427.5
$ws.Range("J113").Value = 503.81818
$ws.Range("K113").Value = 1282.5
$ws.Range("L113").Value = 1511.45454
$ws.Range("M113").Value = 887.5
$ws.Range("N113").Value = -5851.45454
$ws.Range("H132").Value = 10471.286
$ws.Range("I132").Value = 15575
$ws.Range("J132").Value = 3666.3333
$ws.Range("K132").Value = 46725
$ws.Range("L132").Value = 10998.9999
$ws.Range("M132").Value = -44195
$ws.Range("N132").Value = -16058.9999
